# Apply translation / formatting improvements to vpc_routerTable.xlsx
# per commit "Update files based on cn190115 and Improve Translation for console"
#
# Only cells whose content actually changes are touched; everything else is
# left alone so the host re-numbers the shared-string table around the
# surviving (untouched) strings exactly like the source edit did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: "No Route Table Data" -> "No Route Table" ---
$ws.Range("C1").Value = "No Route Table"

# --- Row 4: brand new row content ---
# B4: '自定义表'  -> leading quote becomes quotePrefix, trailing quote stays literal text
$ws.Range("B4").Value = "'自定义表'"

# C4: 'Customized  Route Table'  with "Route" highlighted in red
$ws.Range("C4").Value = "'Customized  Route Table'"
$ws.Range("C4").Characters(13, 5).Font.Color = 255

# D4: empty cell, but carries the same quote-prefix style as B4 / C4
$ws.Range("D4").Value = "'"
$ws.Range("D4").Value = ""

# E4: empty cell with a plain (non-red) Arial font, distinct from the default style
$ws.Range("E4").Font.Name = "Arial"
$ws.Range("E4").Font.Size = 12

# F4: dialog confirm-button label, red
$ws.Range("F4").Value = "确认"
$ws.Range("F4").Font.Color = 255

# --- Row 6: "Associate subnet" -> "Associate Subnet", shown in red ---
$ws.Range("C6").Value = "Associate Subnet"
$ws.Range("C6").Font.Color = 255
$ws.Range("D6").Font.Color = 255

# --- Row 8: "Begin to get the route table list under vpc" -> new wording, red + quote-prefix ---
$ws.Range("C8").Value = "'Begin to get the Route Table Lists of the VPC'"
$ws.Range("C8").Font.Color = 255
$ws.Range("D8").Value = "'"
$ws.Range("D8").Value = ""
$ws.Range("D8").Font.Color = 255
$ws.Range("E8").WrapText = $true

# --- Row 9: "The route table list data under vpc are" -> new wording, red + quote-prefix ---
$ws.Range("C9").Value = "'The Route Table Lists of the VPC are"
$ws.Range("C9").Font.Color = 255
$ws.Range("D9").Value = "'"
$ws.Range("D9").Value = ""
$ws.Range("D9").Font.Color = 255
$ws.Range("E9").Font.Name = "Arial"
$ws.Range("E9").Font.Size = 12
$ws.Range("F9").Value = "确认"
$ws.Range("F9").Font.Color = 255

# --- Column widths (closest values reachable through the host's column-width rounding) ---
$ws.Columns.Item(3).ColumnWidth = 36.142857142857146
$ws.Columns.Item(4).ColumnWidth = 38

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / active cell, set last so it becomes the saved cursor position ---
$ws.Range("C18").Select()
